$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Test Case numbers) ---
$ws.Range("A1").Value = "Test Case"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# --- Column B (new "Description" column) ---
$ws.Range("B1").Value = "Description"
$ws.Range("B2").Value = "passing 1 and 0"
$ws.Range("B4").Value = "passing string"
$ws.Range("B3").Value = "passing b/w 2-9"
$ws.Range("B5").Value = "passing some characters"

# --- Column C (Input, shifted from old column B) ---
$ws.Range("C1").Value = "Input"
$ws.Range("C2").Value = 210
$ws.Range("C3").Value = 23
$ws.Range("C4").Value = " sdd"
$ws.Range("C5").Value = "#$"

# --- Column D (new "Expected Output" column) ---
$ws.Range("D1").Value = "Expected Output"
$ws.Range("D2").Value = "a b c"
$ws.Range("D3").Value = "ad bd cd ae be ce af bf cf"
$ws.Range("D4").Value = """"""
$ws.Range("D5").Value = """"""

# --- Column E (new "Actual Output" column) ---
$ws.Range("E1").Value = "Actual Output"
$ws.Range("E2").Value = "enter string 2-9"
$ws.Range("E3").Value = "ad bd cd ae be ce af bf cf"
$ws.Range("E4").Value = """"""
$ws.Range("E5").Value = """"""

# --- Column F (Result, shifted from old column D) ---
$ws.Range("F1").Value = "Result"
$ws.Range("F2").Value = "FAIL"
$ws.Range("F3").Value = "PASS"
$ws.Range("F4").Value = "PASS"
$ws.Range("F5").Value = "PASS"

# --- Row 6 (new row, appended last) ---
$ws.Range("B6").Value = "with negative input"
$ws.Range("C6").Value = -23
$ws.Range("D6").Value = "empty list"
$ws.Range("E6").Value = "[] "
$ws.Range("F6").Value = "PASS"

# --- Column widths (best achievable given COM rounding precision) ---
$ws.Columns("B").ColumnWidth = 23.333333333333332
$ws.Columns("C").ColumnWidth = 6
$ws.Columns("D").ColumnWidth = 22
$ws.Columns("E").ColumnWidth = 22.666666666666668

# --- Selection ---
$ws.Range("F7").Select()
